$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Recommandations")

# --- Update 'Recommandations' sheet ---
$ws1.Range("D2").Value = 3376.27
$ws1.Range("E2").Value = 108.97
$ws1.Range("D3").Value = 2835
$ws1.Range("E3").Value = 700
$ws1.Range("D4").Value = 2820
$ws1.Range("E4").Value = 715
$ws1.Range("D5").Value = 2467.48
$ws1.Range("E5").Value = 611.62
$ws1.Range("D6").Value = 2028.32
$ws1.Range("E6").Value = 501.93
$ws1.Range("D7").Value = 1491.85
$ws1.Range("E7").Value = 367.37
$ws1.Range("D8").Value = 1439.26
$ws1.Range("E8").Value = 364.71
$ws1.Range("D9").Value = 717.48
$ws1.Range("E9").Value = 180.08
$ws1.Range("D10").Value = 562.07
$ws1.Range("E10").Value = 140.85
$ws1.Range("D11").Value = 561.63
$ws1.Range("E11").Value = 140.92
$ws1.Range("D12").Value = 552.39
$ws1.Range("E12").Value = 138.42
$ws1.Range("D13").Value = 517.64
$ws1.Range("E13").Value = 128.88
$ws1.Range("D14").Value = 439.84
$ws1.Range("E14").Value = 107.08
$ws1.Range("D16").Value = 383.75
$ws1.Range("E16").Value = 97.25
$ws1.Range("A18").Value = 'BRVM - INDUSTRIE              (**)'
$ws1.Range("D18").Value = 219.1
$ws1.Range("E18").Value = 219.1
$ws1.Range("A19").Value = 'BRVM - INDUSTRIE                  (**)'
$ws1.Range("D19").Value = 218.66
$ws1.Range("E19").Value = 218.66
$ws1.Range("A20").Value = 'BRVM-PRINCIPAL                 (**)'
$ws1.Range("D20").Value = 204.98
$ws1.Range("E20").Value = 204.98
$ws1.Range("A23").Value = 'BRVM - CONSOMMATION DE BASE           (**)'
$ws1.Range("D23").Value = 193.79
$ws1.Range("E23").Value = 193.79
$ws1.Range("B24").Value = 3
$ws1.Range("D24").Value = 18.48
$ws1.Range("A25").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Range("B25").Value = 3
$ws1.Range("D25").Value = 9.43
$ws1.Range("E25").Value = 4.08
$ws1.Range("F25").Value = '🟢 Achat'
$ws1.Range("G25").Value = '✅ Renforcer'
$ws1.Range("A27").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("D27").Value = 7.21
$ws1.Range("E27").Value = 7.21
$ws1.Range("A28").Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Range("B28").Value = 1
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = 5.01
$ws1.Range("E28").Value = -1.99
$ws1.Range("G28").Value = '👀 À surveiller'
$ws1.Range("A30").Value = 'SICOR CI (SICC)'
$ws1.Range("C30").Value = 1
$ws1.Range("D30").Value = 4.68
$ws1.Range("E30").Value = -2.78
$ws1.Range("G30").Value = '👀 À surveiller'
$ws1.Range("A31").Value = 'CIE CI (CIEC)'
$ws1.Range("D31").Value = 4.54
$ws1.Range("E31").Value = 4.54
$ws1.Range("A32").Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Range("D32").Value = 3.47
$ws1.Range("E32").Value = 3.47
$ws1.Range("A33").Value = 'SOGB CI (SOGC)'
$ws1.Range("C33").Value = 0
$ws1.Range("D33").Value = 2.07
$ws1.Range("E33").Value = 2.07
$ws1.Range("G33").Value = '➖ Neutre'
$ws1.Range("A34").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 0
$ws1.Range("D34").Value = 1.45
$ws1.Range("E34").Value = 1.45
$ws1.Range("G34").Value = '➖ Neutre'
$ws1.Range("A35").Value = 'SAPH CI (SPHC)'
$ws1.Range("B35").Value = 1
$ws1.Range("C35").Value = 0
$ws1.Range("D35").Value = 1.44
$ws1.Range("E35").Value = 1.44
$ws1.Range("A36").Value = 'SAFCA CI (SAFC)'
$ws1.Range("B36").Value = 2
$ws1.Range("C36").Value = 2
$ws1.Range("D36").Value = 0.09
$ws1.Range("E36").Value = 7.44
$ws1.Range("A37").Value = 'TOTAL'
$ws1.Range("C37").Value = 4
$ws1.Range("D37").Value = 0
$ws1.Range("E37").Value = 0
$ws1.Range("A38").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = -0.65
$ws1.Range("E38").Value = 5.49
$ws1.Range("A39").Value = 'PALM CI (PALC)'
$ws1.Range("D39").Value = -1.99
$ws1.Range("E39").Value = -1.99
$ws1.Range("A40").Value = 'BANK OF AFRICA ML (BOAM)'
$ws1.Range("D40").Value = -2.19
$ws1.Range("E40").Value = -2.19
$ws1.Range("A41").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("D41").Value = -2.5
$ws1.Range("E41").Value = -2.5
$ws1.Range("A42").Value = 'AIR LIQUIDE CI (SIVC)'
$ws1.Range("D42").Value = -3.5
$ws1.Range("E42").Value = -3.5
$ws1.Range("A43").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("D43").Value = -3.85
$ws1.Range("E43").Value = -3.85
$ws1.Range("A44").Value = 'TOTALENERGIES MARKETING CI (TTLC)'
$ws1.Range("D44").Value = -7.35
$ws1.Range("E44").Value = -7.35
$ws1.Range("A45").Value = 'SICABLE CI (CABC)'
$ws1.Range("D45").Value = -10.02
$ws1.Range("E45").Value = -7.37
$ws1.Range("G45").Value = '👀 À surveiller'
$ws1.Range("A46").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B46").Value = 0
$ws1.Range("C46").Value = 2
$ws1.Range("D46").Value = -10.12
$ws1.Range("E46").Value = -4.46
$ws1.Range("F46").Value = '🟡 Observer'
$ws1.Range("G46").Value = '➖ Neutre'
$ws1.Range("A47").Value = 'SETAO CI (STAC)'
$ws1.Range("D47").Value = -10.33
$ws1.Range("E47").Value = -7.08
$ws1.Range("A48").Value = 'CFAO MOTORS CI (CFAC)'
$ws1.Range("B48").Value = 1
$ws1.Range("C48").Value = 3
$ws1.Range("D48").Value = -13.09
$ws1.Range("E48").Value = 7.46
$ws1.Range("F48").Value = '🔴 Vente'
$ws1.Range("G48").Value = '⚠️ Risque de décrochage'

$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Update 'Top_YTD' sheet ---
$ws2.Range("B2").Value = 9158945.78
$ws2.Range("B3").Value = 427612
$ws2.Range("B4").Value = 419462
$ws2.Range("B5").Value = 263977.93
$ws2.Range("B6").Value = 135594.68
$ws2.Range("B7").Value = 49932.56
$ws2.Range("B8").Value = 44596.24
$ws2.Range("B9").Value = 5987.89
$ws2.Range("B10").Value = 3246.4
$ws2.Range("B11").Value = 3240.21
